$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 10.03
$ws.Range("D3").Value = 13.81
$ws.Range("D4").Value = 17.11
$ws.Range("D5").Value = 9.81
$ws.Range("D6").Value = 8.75
$ws.Range("D7").Value = 0.61
$ws.Range("D8").Value = 0.93
$ws.Range("D9").Value = 1.66
$ws.Range("D10").Value = 1.27
$ws.Range("D11").Value = 3.45
$ws.Range("D12").Value = 2.85
$ws.Range("D13").Value = 0.71
$ws.Range("D14").Value = 0.22
$ws.Range("D15").Value = 0.63
$ws.Range("D17").Value = 0.11
$ws.Range("D18").Value = 0.68
$ws.Range("D19").Value = 7.82
$ws.Range("D20").Value = 2.75
$ws.Range("D23").Value = 3.51
$ws.Range("D24").Value = 30.94
$ws.Range("D25").Value = 6.57
$ws.Range("D26").Value = 19.61
$ws.Range("D27").Value = 2.23
$ws.Range("D28").Value = 0.76
$ws.Range("D29").Value = 4.54
$ws.Range("D30").Value = 3.98
$ws.Range("D33").Value = 0.83
$ws.Range("D34").Value = 1.02
$ws.Range("D35").Value = 15.16
$ws.Range("D36").Value = 4.83
$ws.Range("D37").Value = 6.64
$ws.Range("D38").Value = 8.119999999999999
$ws.Range("D39").Value = 4.39
$ws.Range("D40").Value = 0.76
$ws.Range("D41").Value = 1.18
$ws.Range("D42").Value = 55.18
$ws.Range("D43").Value = 0.95
$ws.Range("D44").Value = 0.95
$ws.Range("D45").Value = 4.76
$ws.Range("D46").Value = 1.05
$ws.Range("D47").Value = 3.36
$ws.Range("D48").Value = 6.35
$ws.Range("D49").Value = 1.53
$ws.Range("D51").Value = 4.84
$ws.Range("D52").Value = 0.9399999999999999
$ws.Range("D53").Value = 4.36
$ws.Range("D54").Value = 3.36
$ws.Range("D55").Value = 1.26
$ws.Range("D56").Value = 1.54
$ws.Range("D57").Value = 1.41
$ws.Range("D58").Value = 1.41
$ws.Range("D59").Value = 1.7
$ws.Range("D60").Value = 1.63
$ws.Range("D61").Value = 2.61
$ws.Range("D63").Value = 1.87
$ws.Range("D64").Value = 5.27
$ws.Range("D65").Value = 28.32
$ws.Range("D66").Value = 9.789999999999999
$ws.Range("D67").Value = 1.55
$ws.Range("D68").Value = 1.04
$ws.Range("D69").Value = 3.52
$ws.Range("D70").Value = 4.17
$ws.Range("D71").Value = 0.6
$ws.Range("D72").Value = 0.83
$ws.Range("D73").Value = 2.55
$ws.Range("D74").Value = 2.82
$ws.Range("D75").Value = 6.94
$ws.Range("D76").Value = 4.9
$ws.Range("D77").Value = 0.96
$ws.Range("D78").Value = 1.13
$ws.Range("D79").Value = 6.06
$ws.Range("D80").Value = 3.87
$ws.Range("D81").Value = 3.48
$ws.Range("D82").Value = 1.62
$ws.Range("D83").Value = 1.37
$ws.Range("D84").Value = 10.09
$ws.Range("D85").Value = 7.42
$ws.Range("D86").Value = 17.85
$ws.Range("D87").Value = 1.15
$ws.Range("D88").Value = 9.539999999999999
$ws.Range("D89").Value = 9.130000000000001
$ws.Range("D90").Value = 3.08
$ws.Range("D91").Value = 1.99
$ws.Range("D92").Value = 0.66
$ws.Range("D93").Value = 5.31
$ws.Range("D95").Value = 0.76
$ws.Range("D96").Value = 0.5600000000000001
$ws.Range("D97").Value = 1.77
$ws.Range("D98").Value = 0.86
$ws.Range("D99").Value = 0.64
$ws.Range("D100").Value = 1.25
